$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Results1")
$ws1.Range("E2").Value = 0.01058877300238237
$ws1.Range("E3").Value = 0.01084896300744731
$ws1.Range("E4").Value = 0.0112079660029849
$ws1.Range("E5").Value = 0.01250027200148907
$ws1.Range("E6").Value = 0.01222316898929421
$ws1.Range("E7").Value = 0.08242858799349051
$ws1.Range("E8").Value = 0.04035628600104246
$ws1.Range("E9").Value = 0.03688739400240593
$ws1.Range("E10").Value = 0.05002444199635647
$ws1.Range("E11").Value = 0.04121896599826869
$ws1.Range("E12").Value = 0.2398643740016269
$ws1.Range("E13").Value = 0.2510825019999174
$ws1.Range("E14").Value = 0.2618665240006521
$ws1.Range("E15").Value = 0.2779105600056937
$ws1.Range("E16").Value = 0.2711487049964489
$ws1.Range("E17").Value = 0.00736568401043769
$ws1.Range("E18").Value = 0.004415485993376933
$ws1.Range("E19").Value = 0.005001676006941125
$ws1.Range("E20").Value = 0.005535774995223619
$ws1.Range("E21").Value = 0.005772446005721577
$ws1.Range("E22").Value = 0.1415659080084879
$ws1.Range("E23").Value = 0.06523129300330766
$ws1.Range("E24").Value = 0.03888741899572778
$ws1.Range("E25").Value = 0.08796854900720064
$ws1.Range("E26").Value = 0.0454000379977515
$ws1.Range("E27").Value = 0.3342516929988051
$ws1.Range("E28").Value = 0.3038843689864734
$ws1.Range("E29").Value = 0.316666015991359
$ws1.Range("E30").Value = 0.3439362789940787
$ws1.Range("E31").Value = 0.3340900580078596
$ws1.Range("E32").Value = 6.546001648996025
$ws1.Range("E33").Value = 6.164249348003068
$ws1.Range("E34").Value = 6.941108212995459
$ws1.Range("E35").Value = 8.295261053004651
$ws1.Range("E36").Value = 8.304069053003332
$ws1.Range("E37").Value = 53.78559617900464
$ws1.Range("E38").Value = 16.87291693800944
$ws1.Range("E39").Value = 11.89524010100286
$ws1.Range("E40").Value = 28.03955891799706
$ws1.Range("E41").Value = 9.322999898009584
$ws1.Range("E42").Value = 217.7611038970063
$ws1.Range("E43").Value = 38.38030440399598
$ws1.Range("E44").Value = 25.96025929099414
$ws1.Range("E45").Value = 62.92718287500611
$ws1.Range("E46").Value = 33.84845888100972
$ws1.Range("E47").Value = 11.67326548800338
$ws1.Range("E48").Value = 11.33017684199149
$ws1.Range("E49").Value = 12.00220040499698
$ws1.Range("E50").Value = 14.11401870299596
$ws1.Range("E51").Value = 13.49605374099337
$ws1.Range("E52").Value = 45.40279168200505
$ws1.Range("E53").Value = 15.78876889799722
$ws1.Range("E54").Value = 12.22609734999423
$ws1.Range("E55").Value = 27.67630504899716
$ws1.Range("E56").Value = 18.08479998000257
$ws1.Range("E57").Value = 209.0494304309977
$ws1.Range("E58").Value = 49.21516536599665
$ws1.Range("E59").Value = 46.52588223299244
$ws1.Range("E60").Value = 87.71839000900218
$ws1.Range("E61").Value = 65.39027614500083

$ws2 = $wb.Worksheets.Item("Results2")
$ws2.Range("E2").Value = 0.4342660840047756
$ws2.Range("E3").Value = 0.6151491339987842
$ws2.Range("E4").Value = 0.4391592110041529
$ws2.Range("E5").Value = 0.7068641749938251
$ws2.Range("E6").Value = 0.4790072169998894
$ws2.Range("E7").Value = 37.39261812900077
$ws2.Range("E8").Value = 35.25426222800161
$ws2.Range("E9").Value = 26.91546964700683
$ws2.Range("E10").Value = 49.25394255800347
$ws2.Range("E11").Value = 34.30945669599168
$ws2.Range("E12").Value = 140.7796935850056
$ws2.Range("E13").Value = 72.99056495100376
$ws2.Range("E14").Value = 34.555685653002
$ws2.Range("E15").Value = 104.8992867180059
$ws2.Range("E16").Value = 43.08881236200978
$ws2.Range("E17").Value = 0.3835997210117057
$ws2.Range("E18").Value = 0.530557594000129
$ws2.Range("E19").Value = 0.5115735040017171
$ws2.Range("E20").Value = 0.5822852819983382
$ws2.Range("E21").Value = 0.5646136700088391
$ws2.Range("E22").Value = 18.74061004299438
$ws2.Range("E23").Value = 15.22049036099634
$ws2.Range("E24").Value = 14.97122426598798
$ws2.Range("E25").Value = 20.19500284499372
$ws2.Range("E26").Value = 17.82521583899506
$ws2.Range("E27").Value = 99.75754632000462
$ws2.Range("E28").Value = 34.89610663799976
$ws2.Range("E29").Value = 29.79970340800355
$ws2.Range("E30").Value = 49.33988917199895
$ws2.Range("E31").Value = 34.20889212700422
$ws2.Range("E32").Value = 1.11903757500113
$ws2.Range("E33").Value = 1.252294510006323
$ws2.Range("E34").Value = 1.179804424013128
$ws2.Range("E35").Value = 1.186803327000234
$ws2.Range("E36").Value = 1.197759600006975
$ws2.Range("E37").Value = 22.33944455900928
$ws2.Range("E38").Value = 19.29329607699765
$ws2.Range("E39").Value = 14.09847174200695
$ws2.Range("E40").Value = 23.9155814639962
$ws2.Range("E41").Value = 16.03574053799093
$ws2.Range("E42").Value = 108.8301497609937
$ws2.Range("E43").Value = 60.61442365699622
$ws2.Range("E44").Value = 56.19837358299992
$ws2.Range("E45").Value = 79.34198274499795
$ws2.Range("E46").Value = 60.95489245699719
$ws2.Range("E47").Value = 2.677357029999257
$ws2.Range("E48").Value = 3.536816116000409
$ws2.Range("E49").Value = 3.836212934009382
$ws2.Range("E50").Value = 4.111728917996516
$ws2.Range("E51").Value = 4.265971529006492
$ws2.Range("E52").Value = 16.28540998599783
$ws2.Range("E53").Value = 35.54263880700455
$ws2.Range("E54").Value = 11.87593407300301
$ws2.Range("E55").Value = 43.9188199229975
$ws2.Range("E56").Value = 13.09701940500236
$ws2.Range("E57").Value = 62.34008274899679
$ws2.Range("E58").Value = 46.07152797299204
$ws2.Range("E59").Value = 14.36279468800058
$ws2.Range("E60").Value = 86.04831000098784
$ws2.Range("E61").Value = 22.92530938600248
$ws2.Range("E62").Value = 3.583613772003446
$ws2.Range("E63").Value = 4.561874115999672
$ws2.Range("E64").Value = 5.024546314001782
$ws2.Range("E65").Value = 5.365207245995407
$ws2.Range("E66").Value = 5.57607519600424
$ws2.Range("E67").Value = 21.76227430999279
$ws2.Range("E68").Value = 11.36937612800102
$ws2.Range("E69").Value = 9.444409079995239
$ws2.Range("E70").Value = 13.306516851997
$ws2.Range("E71").Value = 9.430045143002644
$ws2.Range("E72").Value = 24.09288049700262
$ws2.Range("E73").Value = 25.20232958500856
$ws2.Range("E74").Value = 21.93620225100312
$ws2.Range("E75").Value = 29.49400094999874
$ws2.Range("E76").Value = 19.41906819101132
$ws2.Range("E77").Value = 2.26145975801046
$ws2.Range("E78").Value = 3.814742636997835
$ws2.Range("E79").Value = 4.151694428990595
$ws2.Range("E80").Value = 4.355832175991964
$ws2.Range("E81").Value = 4.438339175991132
$ws2.Range("E82").Value = 10.12095094200049
$ws2.Range("E83").Value = 10.73493182400125
$ws2.Range("E84").Value = 10.31960234799772
$ws2.Range("E85").Value = 12.00344322899764
$ws2.Range("E86").Value = 9.669544940988999
$ws2.Range("E87").Value = 21.88447783500305
$ws2.Range("E88").Value = 32.28776554900105
$ws2.Range("E89").Value = 31.61455654499878
$ws2.Range("E90").Value = 35.67640119600401
$ws2.Range("E91").Value = 20.67004490000545
